$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 130, shifting existing rows 130-137 down to 131-138
$ws.Rows(130).Insert()

# Populate the new row 130 with the new weekly price entry
$ws.Range("A130").Value = 11
$ws.Range("B130").Value = "Vega Monumental Concepción"
$ws.Range("C130").Value = "Bíobío"
$ws.Range("D130").Value = 44931
$ws.Range("E130").Value = 8
$ws.Range("F130").Value = 100112028
$ws.Range("G130").Value = "Sandia"
$ws.Range("H130").Value = "Sin especificar"
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 2500
$ws.Range("K130").Value = 1800
$ws.Range("L130").Value = 2000
$ws.Range("M130").Value = 1920
$ws.Range("N130").Value = "$/unidad"
$ws.Range("O130").Value = "Región de O'Higgins"
$ws.Range("P130").Value = 1920
$ws.Range("Q130").Value = 1
$ws.Range("R130").Value = "Hortaliza"
